# Append a new data row (row 3) to Sheet1, mirroring the existing
# MIGRATION DATE / FINANCIAL INSTITUTION NAME / ENTITY ID / ADDRESS columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a date-shaped string ("2025-10-13") that, like the existing
# row 2 value ("2025-10-17"), must stay literal text rather than being
# auto-converted into a date serial number. Temporarily force a text number
# format before assigning the value, then clear the formatting again so the
# cell ends up with the same (default) style as its neighbours.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-10-13"
$ws.Range("A3").ClearFormats()

$ws.Range("B3").Value = "xxss"
$ws.Range("C3").Value = "456CDX012"
$ws.Range("D3").Value = "nana Nagar"
